$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "69.434.28"
$ws.Cells.Item(2, 5).Value = "  -2.80%  "
$ws.Cells.Item(3, 4).Value = "3.540.37"
$ws.Cells.Item(3, 5).Value = "  -2.51%  "
$ws.Cells.Item(4, 5).Value = "  +0.12%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "578.08"
$ws.Cells.Item(5, 5).Value = "  -1.75%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "171.14"
$ws.Cells.Item(6, 5).Value = "  -4.83%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.620"
$ws.Cells.Item(7, 5).Value = "  +0.95%  "
$ws.Cells.Item(8, 4).Value = "3.538.42"
$ws.Cells.Item(8, 5).Value = "  -2.38%  "
$ws.Cells.Item(9, 5).Value = "  +0.16%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.190"
$ws.Cells.Item(10, 5).Value = "  -5.96%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "6.52"
$ws.Cells.Item(11, 5).Value = "  +11.28%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.600"
$ws.Cells.Item(12, 5).Value = "  -0.91%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "47.24"
$ws.Cells.Item(13, 5).Value = "  -4.59%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.0000275"
$ws.Cells.Item(14, 5).Value = "  -3.53%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "686.82"
$ws.Cells.Item(15, 5).Value = "  +0.67%  "
$ws.Cells.Item(16, 4).Value = "4.123.90"
$ws.Cells.Item(16, 5).Value = "  -2.34%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "8.79"
$ws.Cells.Item(17, 5).Value = "  -2.02%  "
$ws.Cells.Item(18, 2).Value = "WrappedBTC"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(18, 4).Value = "69.633.73"
$ws.Cells.Item(18, 5).Value = "  -2.84%  "
$ws.Cells.Item(19, 2).Value = "WrappedEther"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(19, 4).Value = "3.566.14"
$ws.Cells.Item(19, 5).Value = "  -1.80%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.121"
$ws.Cells.Item(20, 5).Value = "  -1.17%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "17.36"
$ws.Cells.Item(21, 5).Value = "  -4.88%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "11.14"
$ws.Cells.Item(22, 5).Value = "  -3.91%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.914"
$ws.Cells.Item(23, 5).Value = "  -1.93%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "16.65"
$ws.Cells.Item(24, 5).Value = "  -5.98%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "97.78"
$ws.Cells.Item(25, 5).Value = "  -5.13%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.84"
$ws.Cells.Item(26, 5).Value = "  -4.21%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.69"
$ws.Cells.Item(27, 5).Value = "  -4.85%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.999"
$ws.Cells.Item(28, 5).Value = "  -0.07%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.44"
$ws.Cells.Item(29, 5).Value = "  -5.14%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "33.67"
$ws.Cells.Item(30, 5).Value = "  -3.74%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "8.91"
$ws.Cells.Item(31, 5).Value = "  -2.97%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.18"
$ws.Cells.Item(32, 5).Value = "  -5.87%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "7.30"
$ws.Cells.Item(33, 5).Value = "  +0.94%  "
$ws.Cells.Item(34, 5).Value = "  -5.59%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.84"
$ws.Cells.Item(35, 5).Value = "  -6.79%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "566.36"
$ws.Cells.Item(36, 5).Value = "  -1.20%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "10.85"
$ws.Cells.Item(37, 5).Value = "  -4.09%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.105"
$ws.Cells.Item(38, 5).Value = "  -3.91%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "57.52"
$ws.Cells.Item(39, 5).Value = "  -3.24%  "
$ws.Cells.Item(40, 5).Value = "  -0.06%  "
$ws.Cells.Item(41, 2).Value = "Maker"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(41, 4).Value = "3.493.30"
$ws.Cells.Item(41, 5).Value = "  -4.68%  "
$ws.Cells.Item(42, 2).Value = "VeChain"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0442"
$ws.Cells.Item(42, 5).Value = "  -5.87%  "
$ws.Cells.Item(43, 2).Value = "Kaspa"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.138"
$ws.Cells.Item(43, 5).Value = "  -3.16%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.337"
$ws.Cells.Item(44, 5).Value = "  -2.45%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "33.48"
$ws.Cells.Item(45, 5).Value = "  -5.53%  "
$ws.Cells.Item(46, 4).Value = "0.0₃0707"
$ws.Cells.Item(46, 5).Value = "  -6.65%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.90"
$ws.Cells.Item(47, 5).Value = "  +3.96%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.58"
$ws.Cells.Item(48, 5).Value = "  -5.92%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.133"
$ws.Cells.Item(49, 5).Value = "  -0.17%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "133.80"
$ws.Cells.Item(50, 5).Value = "  +1.96%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.149"
$ws.Cells.Item(51, 5).Value = "  -0.52%  "
